$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update train name in B3 from "Swarna Shatabdi" to "Uhl Janstb Spl"
$ws.Range("B3").Value = "Uhl Janstb Spl"

# Update the active selection to B4 (matches saved view state in the diff)
$ws.Range("B4").Select()
